# Pequenos ajustes e correcoes v2
# Appends 3 new login log rows (maria logged in three more times) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @(
    @("2026-01-06", "17:02:50", "maria", "maria@teste.com"),
    @("2026-01-06", "17:07:19", "maria", "maria@teste.com"),
    @("2026-01-06", "17:07:36", "maria", "maria@teste.com")
)

$startRow = 21
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $rowRange = $ws.Range("A$r`:D$r")

    # Force the date-looking text (e.g. "2026-01-06") to stay plain text
    # instead of letting Excel auto-convert it to a numeric date serial,
    # then drop the formatting override so the cell keeps the default
    # (unstyled) look of the other data rows.
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    $rowRange.ClearFormats()
}
